$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume columns keep their original text formatting so that
# Excel does not auto-convert numeric-looking strings (e.g. "1.004") into numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '26.892.18'
$ws.Range("E2").Value = '  -0.19%  '

$ws.Range("D3").Value = '1.813.82'
$ws.Range("E3").Value = '  +1.49%  '

$ws.Range("D4").Value = '1.004'
$ws.Range("E4").Value = '  -0.68%  '

$ws.Range("D5").Value = '311.32'
$ws.Range("E5").Value = '  -0.89%  '

$ws.Range("E6").Value = '  -0.63%  '

$ws.Range("D7").Value = '0.4292'
$ws.Range("E7").Value = '  +1.35%  '

$ws.Range("D8").Value = '0.3690'
$ws.Range("E8").Value = '  +2.34%  '

$ws.Range("D9").Value = '0.07237'
$ws.Range("E9").Value = '  +1.20%  '

$ws.Range("D10").Value = '0.8613'
$ws.Range("E10").Value = '  +2.39%  '

$ws.Range("B11").Value = 'WrappedEther'
$ws.Range("C11").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D11").Value = '2.015.26'
$ws.Range("E11").Value = '  +8.76%  '

$ws.Range("B12").Value = 'Solana'
$ws.Range("C12").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D12").Value = '21.06'
$ws.Range("E12").Value = '  +3.85%  '

$ws.Range("D13").Value = '6.639'
$ws.Range("E13").Value = '  +4.57%  '

$ws.Range("D14").Value = '5.397'
$ws.Range("E14").Value = '  +2.67%  '

$ws.Range("E15").Value = '  +1.11%  '

$ws.Range("D16").Value = '80.68'
$ws.Range("E16").Value = '  +1.39%  '

$ws.Range("D17").Value = '1.005'
$ws.Range("E17").Value = '  -0.88%  '

$ws.Range("D18").Value = '0.000008926'
$ws.Range("E18").Value = '  +2.80%  '

$ws.Range("E19").Value = '  -0.67%  '

$ws.Range("D20").Value = '15.19'
$ws.Range("E20").Value = '  +1.78%  '

$ws.Range("D21").Value = '26.948.79'
$ws.Range("E21").Value = '  -0.75%  '

$ws.Range("D22").Value = '5.180'
$ws.Range("E22").Value = '  +2.53%  '

$ws.Range("D23").Value = '11.06'
$ws.Range("E23").Value = '  -0.10%  '

$ws.Range("D24").Value = '2.236.17'
$ws.Range("E24").Value = '  +7.88%  '

$ws.Range("D25").Value = '153.70'
$ws.Range("E25").Value = '  +0.29%  '

$ws.Range("D26").Value = '1.882'
$ws.Range("E26").Value = '  -3.25%  '

$ws.Range("D27").Value = '18.23'
$ws.Range("E27").Value = '  +0.21%  '

$ws.Range("D28").Value = '5.209'
$ws.Range("E28").Value = '  +3.94%  '

$ws.Range("D29").Value = '114.95'
$ws.Range("E29").Value = '  +0.25%  '

$ws.Range("D30").Value = '1.867'
$ws.Range("E30").Value = '  +15.32%  '

$ws.Range("D31").Value = '0.08943'
$ws.Range("E31").Value = '  +0.07%  '

$ws.Range("D32").Value = '0.7428'
$ws.Range("E32").Value = '  +3.03%  '

$ws.Range("D33").Value = '1.162'
$ws.Range("E33").Value = '  +7.21%  '

$ws.Range("E34").Value = '  +2.31%  '

$ws.Range("D35").Value = '2.799'
$ws.Range("E35").Value = '  -1.79%  '

$ws.Range("D36").Value = '1.008'
$ws.Range("E36").Value = '  -0.28%  '

$ws.Range("D37").Value = '1.117'
$ws.Range("E37").Value = '  +3.49%  '

$ws.Range("D38").Value = '0.05215'
$ws.Range("E38").Value = '  +2.70%  '

$ws.Range("D39").Value = '0.01922'
$ws.Range("E39").Value = '  +1.33%  '

$ws.Range("D40").Value = '0.5075'
$ws.Range("E40").Value = '  +2.71%  '

$ws.Range("E41").Value = '  +1.70%  '

$ws.Range("D42").Value = '2.721'
$ws.Range("E42").Value = '  +8.44%  '

$ws.Range("E43").Value = '  +7.14%  '

$ws.Range("D44").Value = '8.244'
$ws.Range("E44").Value = '  +3.61%  '

$ws.Range("D45").Value = '106.70'
$ws.Range("E45").Value = '  +2.16%  '

$ws.Range("D46").Value = '10.41'
$ws.Range("E46").Value = '  +2.62%  '

$ws.Range("D47").Value = '1.004'
$ws.Range("E47").Value = '  -0.71%  '

$ws.Range("D48").Value = '1.655'
$ws.Range("E48").Value = '  +5.09%  '

$ws.Range("B49").Value = 'Decentraland'
$ws.Range("C49").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D49").Value = '0.4573'
$ws.Range("E49").Value = '  +2.11%  '

$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D50").Value = '0.06294'
$ws.Range("E50").Value = '  +0.24%  '

$ws.Range("D51").Value = '1.799'
$ws.Range("E51").Value = '  +6.12%  '
